# Updates the cryptos list with refreshed prices / % volume figures.
# Price values that look like plain decimal numbers are prefixed with a
# leading apostrophe so Excel stores them as text (matching the original
# inline-string cell type) instead of silently converting them to numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 : Bitcoin ---
$ws.Range("D2").Value = "90.576.99"
$ws.Range("E2").Value = "  -0.91%  "

# --- Row 3 : Ethereum ---
$ws.Range("D3").Value = "3.105.82"
$ws.Range("E3").Value = "  -0.84%  "

# --- Row 4 : TetherUSD ---
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  +0.32%  "

# --- Row 5 : Solana ---
$ws.Range("D5").Value = "'237.25"
$ws.Range("E5").Value = "  +7.26%  "

# --- Row 6 : BNB ---
$ws.Range("D6").Value = "'625.57"

# --- Row 7 : XRP ---
$ws.Range("D7").Value = "'1.02"
$ws.Range("E7").Value = "  +4.50%  "

# --- Row 8 : Dogecoin ---
$ws.Range("E8").Value = "  -7.28%  "

# --- Row 10 : LidoStakedEther ---
$ws.Range("D10").Value = "3.381.77"
$ws.Range("E10").Value = "  +8.04%  "

# --- Row 11 : Cardano ---
$ws.Range("D11").Value = "'0.712"
$ws.Range("E11").Value = "  -3.19%  "

# --- Row 12 : TRON ---
$ws.Range("E12").Value = "  +3.45%  "

# --- Row 13 : Avalanche ---
$ws.Range("D13").Value = "'36.21"
$ws.Range("E13").Value = "  +2.30%  "

# --- Row 14 : ShibaInu ---
$ws.Range("E14").Value = "  -4.60%  "

# --- Row 15 : Toncoin ---
$ws.Range("D15").Value = "'5.57"
$ws.Range("E15").Value = "  +2.33%  "

# --- Row 16 : WrappedBTC ---
$ws.Range("D16").Value = "90.224.95"
$ws.Range("E16").Value = "  -1.13%  "

# --- Row 17 : WrappedliquidstakedEther2.0 ---
$ws.Range("D17").Value = "3.676.22"
$ws.Range("E17").Value = "  -0.98%  "

# --- Row 18 : WrappedEther ---
$ws.Range("D18").Value = "3.097.52"
$ws.Range("E18").Value = "  -0.95%  "

# --- Row 19 : SuiNetwork ---
$ws.Range("D19").Value = "'3.71"
$ws.Range("E19").Value = "  -2.08%  "

# --- Row 20 : Chainlink ---
$ws.Range("D20").Value = "'14.25"
$ws.Range("E20").Value = "  +0.40%  "

# --- Row 21 : PEPE ---
$ws.Range("D21").Value = "'0.0000212"
$ws.Range("E21").Value = "  -6.29%  "

# --- Row 22 : BitcoinCash ---
$ws.Range("D22").Value = "'447.61"
$ws.Range("E22").Value = "  +1.54%  "

# --- Row 23 & 24 : Polkadot/Uniswap swapped ranks ---
$ws.Range("B23").Value = "Polkadot"
$ws.Range("C23").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D23").Value = "'5.63"
$ws.Range("E23").Value = "  +8.51%  "

$ws.Range("B24").Value = "Uniswap"
$ws.Range("C24").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D24").Value = "'8.89"
$ws.Range("E24").Value = "  +0.86%  "

# --- Row 25 : NEARProtocol ---
$ws.Range("D25").Value = "'6.08"
$ws.Range("E25").Value = "  -0.28%  "

# --- Row 26 : Litecoin ---
$ws.Range("D26").Value = "'89.77"
$ws.Range("E26").Value = "  +4.09%  "

# --- Row 27 : Aptos ---
$ws.Range("D27").Value = "'12.20"
$ws.Range("E27").Value = "  -1.24%  "

# --- Row 28 : WrappedeETH ---
$ws.Range("D28").Value = "3.258.79"
$ws.Range("E28").Value = "  -1.08%  "

# --- Row 30 : InternetComputer(DFINITY) ---
$ws.Range("D30").Value = "'9.28"
$ws.Range("E30").Value = "  +3.06%  "

# --- Row 31 : Cronos ---
$ws.Range("D31").Value = "'0.159"
$ws.Range("E31").Value = "  -6.40%  "

# --- Row 32 : EthereumClassic ---
$ws.Range("D32").Value = "'27.35"
$ws.Range("E32").Value = "  +14.66%  "

# --- Row 33 : Stellar ---
$ws.Range("D33").Value = "'0.196"
$ws.Range("E33").Value = "  +27.16%  "

# --- Row 34 : Kaspa ---
$ws.Range("E34").Value = "  +4.61%  "

# --- Row 35 : dogwifhat ---
$ws.Range("D35").Value = "'3.79"
$ws.Range("E35").Value = "  -1.51%  "

# --- Row 36 : Bittensor ---
$ws.Range("D36").Value = "'506.42"
$ws.Range("E36").Value = "  -4.93%  "

# --- Row 37 : PancakeSwap ---
$ws.Range("E37").Value = "  +2.78%  "

# --- Row 38 : RenderToken ---
$ws.Range("D38").Value = "'7.01"
$ws.Range("E38").Value = "  -2.92%  "

# --- Row 39 : Fetch.AI ---
$ws.Range("D39").Value = "'1.32"
$ws.Range("E39").Value = "  +1.61%  "

# --- Row 40 : PolygonEcosystemToken ---
$ws.Range("E40").Value = "  +9.59%  "

# --- Row 41 : WhiteBITCoin ---
$ws.Range("D41").Value = "'22.18"
$ws.Range("E41").Value = "  -0.56%  "

# --- Row 42 : Hedera ---
$ws.Range("D42").Value = "'0.0853"
$ws.Range("E42").Value = "  +7.80%  "

# --- Row 43 : USDe ---
$ws.Range("E43").Value = "  +0.01%  "

# --- Row 44 : Binance-PegBSC-USD ---
$ws.Range("D44").Value = "'0.734"
$ws.Range("E44").Value = "  -18.35%  "

# --- Row 45 : MantraDAO ---
$ws.Range("D45").Value = "'3.20"
$ws.Range("E45").Value = "  +33.03%  "

# --- Row 47 : ARBITRUM ---
$ws.Range("E47").Value = "  +10.93%  "

# --- Row 48 : Monero ---
$ws.Range("D48").Value = "'148.76"
$ws.Range("E48").Value = "  +1.77%  "

# --- Row 49 & 50 : OKB/Filecoin swapped ranks ---
$ws.Range("B49").Value = "OKB"
$ws.Range("C49").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D49").Value = "'44.93"
$ws.Range("E49").Value = "  +1.86%  "

$ws.Range("B50").Value = "Filecoin"
$ws.Range("C50").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D50").Value = "'4.50"
$ws.Range("E50").Value = "  +6.55%  "

# --- Row 51 : ImmutableX ---
$ws.Range("E51").Value = "  +3.43%  "
